# HKStock: add readIndexData rows for HSCCI (2016-08-09 .. 2016-08-11)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 265 (index 263) ---
$ws.Range("A265").Value = 263
$ws.Range("B265").Value = "HSCCI"
$ws.Range("C265").Formula = "=""2016-08-09"""
$ws.Range("D265").Value = 3821.06
$ws.Range("E265").Value = 3811.75
$ws.Range("F265").Value = 3850.55
$ws.Range("G265").Value = 3845.15

# --- Row 266 (index 264) ---
$ws.Range("A266").Value = 264
$ws.Range("B266").Value = "HSCCI"
$ws.Range("C266").Formula = "=""2016-08-10"""
$ws.Range("D266").Value = 3853.07
$ws.Range("E266").Value = 3819.14
$ws.Range("F266").Value = 3875.04
$ws.Range("G266").Value = 3848.18

# --- Row 267 (index 265) ---
$ws.Range("A267").Value = 265
$ws.Range("B267").Value = "HSCCI"
$ws.Range("C267").Formula = "=""2016-08-11"""
$ws.Range("D267").Value = 3856.78
$ws.Range("E267").Value = 3833.29
$ws.Range("F267").Value = 3870.33
$ws.Range("G267").Value = 3856.32

# Convert the date formulas in column C to static text values (matches the
# original file's plain shared-string cells, no formula / no number format).
$ws.Range("C265:C267").Copy()
$ws.Range("C265:C267").PasteSpecial(-4163)

# Copy column A's bordered/centered style (from the last existing data row)
# onto the new index cells, matching the style used throughout column A.
$ws.Range("A264").Copy()
$ws.Range("A265:A267").PasteSpecial(-4122)
